$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (Changed) date column for all existing data rows
#    from 45205 (2023-10-06) to 45206 (2023-10-07)
$ws.Range("C2:C204").Value = 45206

# 2. Give row 204 an explicit row height (matches the rest of the data rows)
$ws.Rows.Item(204).RowHeight = 15

# 3. Append the new data row (row 205) describing case "A 47980-2023"
$ws.Range("A205").Value = "A 47980-2023"
$ws.Range("B205").Value = 45204
$ws.Range("C205").Value = 45206
$ws.Range("D205").Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Range("E205").Value = "BOLLEBYGD"
$ws.Range("G205").Value = 4.4
$ws.Range("H205").Value = 0
$ws.Range("I205").Value = 0
$ws.Range("J205").Value = 0
$ws.Range("K205").Value = 0
$ws.Range("L205").Value = 0
$ws.Range("M205").Value = 0
$ws.Range("N205").Value = 0
$ws.Range("O205").Value = 0
$ws.Range("P205").Value = 0
$ws.Range("Q205").Value = 0

# Match number formats used by the other date columns (B/C use style index 1)
$ws.Range("B205:C205").NumberFormat = $ws.Range("B204:C204").NumberFormat

# Column R on data rows uses wrapped text styling even when empty (style index 2)
$ws.Range("R205").WrapText = $true
